$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header order (columns A-P)
$headers = @(
    "Inhaber",
    "total_km_durchgangsstrasse",
    "kb_befreit",
    "AmpelcodePers1",
    "AmpelcodePers2",
    "AmpelcodeOFG1",
    "AmpelcodeOFG2",
    "AmpelcodeOFG3",
    "AmpelcodeGW1",
    "AmpelcodeGW3",
    "AmpelcodeGW4",
    "AmpelcodePers3",
    "AmpelcodePers4",
    "AmpelcodePers5",
    "AmpelcodeOFG5",
    "AmpelcodeGW5"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Data rows keyed by Inhaber name, values keyed by header name.
# Rows are written in the new order: BL, LU, TG (BL first row of data now).
$data = @(
    @{
        "Inhaber" = "BL"
        "total_km_durchgangsstrasse" = 176.2151
        "kb_befreit" = 0
        "AmpelcodePers1" = 162.7862
        "AmpelcodePers2" = 13.4289
        "AmpelcodeOFG1" = 131.8565
        "AmpelcodeOFG2" = 26.0675
        "AmpelcodeOFG3" = 18.2911
        "AmpelcodeGW1" = 169.2787
        "AmpelcodeGW3" = 2.9284
        "AmpelcodeGW4" = 4.008
        "AmpelcodePers3" = $null
        "AmpelcodePers4" = $null
        "AmpelcodePers5" = $null
        "AmpelcodeOFG5" = $null
        "AmpelcodeGW5" = $null
    },
    @{
        "Inhaber" = "LU"
        "total_km_durchgangsstrasse" = 365.898
        "kb_befreit" = 0
        "AmpelcodePers1" = 346.457
        "AmpelcodePers2" = 8.952
        "AmpelcodeOFG1" = 343.37
        "AmpelcodeOFG2" = 8.132999999999999
        "AmpelcodeOFG3" = 14.39
        "AmpelcodeGW1" = 344.948
        "AmpelcodeGW3" = 4.605
        "AmpelcodeGW4" = 16.34
        "AmpelcodePers3" = 7.513
        "AmpelcodePers4" = 2.971
        "AmpelcodePers5" = 0.005
        "AmpelcodeOFG5" = 0.005
        "AmpelcodeGW5" = 0.005
    },
    @{
        "Inhaber" = "TG"
        "total_km_durchgangsstrasse" = 344.233
        "kb_befreit" = 0
        "AmpelcodePers1" = 333.885
        "AmpelcodePers2" = 6.685
        "AmpelcodeOFG1" = 335.28
        "AmpelcodeOFG2" = 5.345
        "AmpelcodeOFG3" = 2.145
        "AmpelcodeGW1" = 284.425
        "AmpelcodeGW3" = 2.928
        "AmpelcodeGW4" = 55.417
        "AmpelcodePers3" = 1.66
        "AmpelcodePers4" = 0.54
        "AmpelcodePers5" = 1.463
        "AmpelcodeOFG5" = 1.463
        "AmpelcodeGW5" = 1.463
    }
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $rowData = $data[$r]
    for ($c = 0; $c -lt $headers.Length; $c++) {
        $colName = $headers[$c]
        $value = $rowData[$colName]
        $cell = $ws.Cells.Item($rowNum, $c + 1)
        if ($null -eq $value) {
            $cell.Value = $null
        } else {
            $cell.Value = $value
        }
    }
}
